$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 621.6667
$ws.Range("I9").Value = 779.8889
$ws.Range("J9").Value = 147
$ws.Range("K9").Value = 779.8889
$ws.Range("L9").Value = 147
$ws.Range("M9").Value = -610.8889
$ws.Range("N9").Value = -485
$ws.Range("H17").Value = 1982.3334
$ws.Range("J17").Value = 1982.3334
$ws.Range("L17").Value = 5947.0002
$ws.Range("N17").Value = -6283.0002
$ws.Range("H51").Value = 11944
$ws.Range("J51").Value = 12599.4
$ws.Range("L51").Value = 12599.4
$ws.Range("N51").Value = -13567.4
$ws.Range("H92").Value = 920.625
$ws.Range("I92").Value = 883
$ws.Range("K92").Value = 883
$ws.Range("M92").Value = 365
$ws.Range("H135").Value = 1929.7778
$ws.Range("I135").Value = 1796
$ws.Range("K135").Value = 16164
$ws.Range("M135").Value = -13629
$ws.Range("H137").Value = 2747.125
$ws.Range("I137").Value = 2711
$ws.Range("K137").Value = 8133
$ws.Range("M137").Value = -5583
$ws.Range("H138").Value = 8362.190000000001
$ws.Range("I138").Value = 1498.8334
$ws.Range("J138").Value = 11107.533
$ws.Range("K138").Value = 4496.5002
$ws.Range("L138").Value = 33322.599
$ws.Range("M138").Value = 643.4997999999996
$ws.Range("N138").Value = -43602.599

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8502103
$ws.Range("I8").Value = 8502103
$ws.Range("K8").Value = 8502103
$ws.Range("M8").Value = -8501959
$ws.Range("H31").Value = 10156.667
$ws.Range("I31").Value = 10156.667
$ws.Range("K31").Value = 10156.667
$ws.Range("M31").Value = -9862.666999999999
$ws.Range("H132").Value = 1857.6666
$ws.Range("I132").Value = 1857.6666
$ws.Range("K132").Value = 5572.9998
$ws.Range("M132").Value = -3042.9998

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1150.8572
$ws.Range("I94").Value = 1531.6666
$ws.Range("J94").Value = 998.5333000000001
$ws.Range("K94").Value = 1531.6666
$ws.Range("L94").Value = 998.5333000000001
$ws.Range("M94").Value = -1080.6666
$ws.Range("N94").Value = -1900.5333
$ws.Range("H99").Value = 1674.1666
$ws.Range("I99").Value = 681.6667
$ws.Range("K99").Value = 681.6667
$ws.Range("M99").Value = 816.3333
$ws.Range("H134").Value = 1573.0769
$ws.Range("I134").Value = 1573.0769
$ws.Range("K134").Value = 4719.2307
$ws.Range("M134").Value = -2184.2307

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5833551.5
$ws.Range("I6").Value = 7000150.5
$ws.Range("K6").Value = 7000150.5
$ws.Range("M6").Value = -7000037.5
$ws.Range("H16").Value = 1105.5
$ws.Range("I16").Value = 1184
$ws.Range("J16").Value = 713
$ws.Range("K16").Value = 1184
$ws.Range("L16").Value = 713
$ws.Range("M16").Value = -897
$ws.Range("N16").Value = -1287
$ws.Range("H22").Value = 4444996.5
$ws.Range("J22").Value = 13333833
$ws.Range("L22").Value = 13333833
$ws.Range("N22").Value = -13334533
$ws.Range("H31").Value = 3002
$ws.Range("I31").Value = 1300
$ws.Range("K31").Value = 1300
$ws.Range("M31").Value = -1005
$ws.Range("H34").Value = 3002
$ws.Range("I34").Value = 1300
$ws.Range("K34").Value = 1300
$ws.Range("M34").Value = -1098
$ws.Range("H86").Value = 4099.6
$ws.Range("I86").Value = 3874.5
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3874.5
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2751.5
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4099.6
$ws.Range("I89").Value = 3874.5
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 19372.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -13756.5
$ws.Range("N89").Value = -36232
$ws.Range("H92").Value = 31804.75
$ws.Range("J92").Value = 31804.75
$ws.Range("L92").Value = 31804.75
$ws.Range("N92").Value = -36796.75
$ws.Range("H113").Value = 1105.5
$ws.Range("I113").Value = 1184
$ws.Range("J113").Value = 713
$ws.Range("K113").Value = 1184
$ws.Range("L113").Value = 713
$ws.Range("M113").Value = 986
$ws.Range("N113").Value = -5053
$ws.Range("H134").Value = 6127.75
$ws.Range("I134").Value = 6127.75
$ws.Range("K134").Value = 18383.25
$ws.Range("M134").Value = -15848.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1828.2
$ws.Range("I5").Value = 1828.2
$ws.Range("K5").Value = 5484.6
$ws.Range("M5").Value = -5372.6
$ws.Range("H39").Value = 6507.7646
$ws.Range("J39").Value = 7288.8667
$ws.Range("L39").Value = 21866.6001
$ws.Range("N39").Value = -22454.6001
$ws.Range("H112").Value = 14499.7
$ws.Range("I112").Value = 7499
$ws.Range("J112").Value = 16249.875
$ws.Range("K112").Value = 22497
$ws.Range("L112").Value = 48749.625
$ws.Range("M112").Value = -21389
$ws.Range("N112").Value = -50965.625
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 1999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17991
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15461
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 18054.889
$ws.Range("I133").Value = 15123.5
$ws.Range("K133").Value = 45370.5
$ws.Range("M133").Value = -40310.5
$ws.Range("H135").Value = 1828.2
$ws.Range("I135").Value = 1828.2
$ws.Range("K135").Value = 16453.8
$ws.Range("M135").Value = -13918.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 705.3333
$ws.Range("I26").Value = 705.3333
$ws.Range("K26").Value = 705.3333
$ws.Range("M26").Value = -410.3333
$ws.Range("H40").Value = 3000.4614
$ws.Range("I40").Value = 2111.3
$ws.Range("K40").Value = 2111.3
$ws.Range("M40").Value = -1975.3
$ws.Range("H62").Value = 22300
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 22300
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 3711.1177
$ws.Range("I132").Value = 3630.5625
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10891.6875
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8361.6875
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 17499.5
$ws.Range("I136").Value = 15000
$ws.Range("J136").Value = 19999
$ws.Range("K136").Value = 45000
$ws.Range("L136").Value = 59997
$ws.Range("M136").Value = -42450
$ws.Range("N136").Value = -65097

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 6860.6
$ws.Range("I55").Value = 625
$ws.Range("J55").Value = 11017.667
$ws.Range("K55").Value = 625
$ws.Range("L55").Value = 11017.667
$ws.Range("M55").Value = -348
$ws.Range("N55").Value = -11571.667
$ws.Range("H132").Value = 3092.7334
$ws.Range("I132").Value = 2730.8462
$ws.Range("K132").Value = 8192.5386
$ws.Range("M132").Value = -5662.5386
$ws.Range("H136").Value = 11513.385
$ws.Range("I136").Value = 11513.385
$ws.Range("L136").Value = 34540.155
$ws.Range("M136").Value = -31990.155
